$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append new row 32 (Testmail #17) ---
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A32").Value = "Ik heb een klacht"
$ws.Range("B32").Value = "mailmind.test@zohomail.eu"
$ws.Range("C32").Value = "Testmail #17: Ik heb een klacht"
$ws.Range("D32").Value = "Klacht / Probleem"

$antwoord = "Beste klant,`r`nBedankt voor het doorsturen van uw klacht. Om uw klacht zo goed mogelijk te kunnen behandelen, ontvangen wij graag meer informatie over de aard van de klacht. Kunt u ons meer details geven over wat er precies is misgegaan? `r`nMet vriendelijke groet,`r`n[Naam] `r`nKlantenservice Team"
$ws.Range("E32").Value = $antwoord

$ws.Range("F32").Value = "2025-06-29 15:17:57"
$ws.Range("G32").Value = "Ja"
$ws.Range("H32").Value = "Nee"
$ws.Range("I32").Value = "Ja"

# Undo the automatic row-height autofit triggered by the multi-line text
# above, so row 32 keeps the default (non-custom) height like every other
# data row in the sheet.
$ws.Rows.Item(32).AutoFit()

# Extend the conditional-formatting ranges so they keep covering the whole
# data range (previously row 2-31, now row 2-32).
$ws.Range("D2:D31").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D32"))
$ws.Range("G2:G31").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G32"))
$ws.Range("H2:H31").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H32"))
$ws.Range("I2:I31").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I32"))

# --- Sheet "Dashboard": bump the "Klacht / Probleem" count from 1 to 2 ---
$ws2 = $wb.Worksheets.Item("Dashboard")
$ws2.Range("B7").Value = 2
